$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 6 (Youtube Salary, Streaming Salary, Salary, Freelance)
$ws.Range("A3:C6").EntireRow.Delete()

# Update row 2 values
$ws.Range("A2").Value = "Payroll"
$ws.Range("B2").Value = 10000
$ws.Range("C2").Value = 45976.333333333336
